# Sync attendance_reports: swap the order of "System" and the email
# in the "Recorded By" column (G) from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System", wherever that exact value occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
